# Applies the "Updated cryptos list" price/volume refresh to Sheet1.
# Column D holds the coin Price as text scraped from the source site
# (sometimes with thousands separators, sometimes plain decimals with
# meaningful trailing zeros); column E holds the Volume(1h) percentage
# text. Both must stay literal text, not be reinterpreted as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price cells (column D) whose new text parses as a plain number ---
# (e.g. "313.55", "0.07282"). Force Text format first so Excel keeps the
# exact scraped digits/trailing zeros instead of coercing to a float.
$priceTextCells = @(
    @{ Cell = 'D5'; Value = '313.55' }
    @{ Cell = 'D7'; Value = '0.5013' }
    @{ Cell = 'D9'; Value = '0.07282' }
    @{ Cell = 'D10'; Value = '0.9082' }
    @{ Cell = 'D12'; Value = '0.07666' }
    @{ Cell = 'D14'; Value = '5.484' }
    @{ Cell = 'D15'; Value = '91.69' }
    @{ Cell = 'D16'; Value = '1.003' }
    @{ Cell = 'D17'; Value = '0.000008722' }
    @{ Cell = 'D20'; Value = '14.52' }
    @{ Cell = 'D22'; Value = '10.84' }
    @{ Cell = 'D23'; Value = '6.601' }
    @{ Cell = 'D24'; Value = '154.33' }
    @{ Cell = 'D25'; Value = '1.879' }
    @{ Cell = 'D26'; Value = '2.236' }
    @{ Cell = 'D28'; Value = '115.39' }
    @{ Cell = 'D29'; Value = '4.911' }
    @{ Cell = 'D30'; Value = '0.08972' }
    @{ Cell = 'D32'; Value = '1.233' }
    @{ Cell = 'D33'; Value = '0.7665' }
    @{ Cell = 'D34'; Value = '4.645' }
    @{ Cell = 'D36'; Value = '2.558' }
    @{ Cell = 'D37'; Value = '0.5610' }
    @{ Cell = 'D39'; Value = '3.014' }
    @{ Cell = 'D40'; Value = '0.05253' }
    @{ Cell = 'D41'; Value = '6.960' }
    @{ Cell = 'D42'; Value = '8.490' }
    @{ Cell = 'D43'; Value = '0.1515' }
    @{ Cell = 'D44'; Value = '111.46' }
    @{ Cell = 'D45'; Value = '10.59' }
    @{ Cell = 'D46'; Value = '0.4817' }
    @{ Cell = 'D49'; Value = '67.42' }
    @{ Cell = 'D50'; Value = '0.06067' }
    @{ Cell = 'D51'; Value = '0.9011' }
)
foreach ($item in $priceTextCells) {
    $rng = $ws.Range($item.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $item.Value
}

# --- Remaining price cells (column D) that already contain extra dots
# (thousands separators) so they can never be parsed as a single number,
# plus every changed Volume(1h) percentage cell (column E). These are safe
# to assign directly. ---
$plainCells = @(
    @{ Cell = 'D2'; Value = '27.876.21' }
    @{ Cell = 'E2'; Value = '  -0.90%  ' }
    @{ Cell = 'D3'; Value = '1.907.73' }
    @{ Cell = 'E3'; Value = '  -0.18%  ' }
    @{ Cell = 'E4'; Value = '  -0.20%  ' }
    @{ Cell = 'E5'; Value = '  -1.00%  ' }
    @{ Cell = 'E6'; Value = '  -0.11%  ' }
    @{ Cell = 'E7'; Value = '  +4.00%  ' }
    @{ Cell = 'E8'; Value = '  -0.13%  ' }
    @{ Cell = 'E9'; Value = '  -1.19%  ' }
    @{ Cell = 'E10'; Value = '  -2.75%  ' }
    @{ Cell = 'E11'; Value = '  -0.08%  ' }
    @{ Cell = 'E12'; Value = '  -1.76%  ' }
    @{ Cell = 'D13'; Value = '1.915.94' }
    @{ Cell = 'E13'; Value = '  +0.47%  ' }
    @{ Cell = 'E14'; Value = '  -0.63%  ' }
    @{ Cell = 'E15'; Value = '  +0.00%  ' }
    @{ Cell = 'E16'; Value = '  -0.21%  ' }
    @{ Cell = 'E17'; Value = '  -1.23%  ' }
    @{ Cell = 'E18'; Value = '  -0.11%  ' }
    @{ Cell = 'D19'; Value = '27.920.14' }
    @{ Cell = 'E19'; Value = '  -0.86%  ' }
    @{ Cell = 'E20'; Value = '  -2.28%  ' }
    @{ Cell = 'E22'; Value = '  -0.76%  ' }
    @{ Cell = 'E23'; Value = '  -0.39%  ' }
    @{ Cell = 'E24'; Value = '  -1.10%  ' }
    @{ Cell = 'E25'; Value = '  -2.20%  ' }
    @{ Cell = 'E26'; Value = '  +5.72%  ' }
    @{ Cell = 'E27'; Value = '  -0.86%  ' }
    @{ Cell = 'E28'; Value = '  -0.98%  ' }
    @{ Cell = 'E29'; Value = '  -1.19%  ' }
    @{ Cell = 'E30'; Value = '  +0.18%  ' }
    @{ Cell = 'E31'; Value = '  -3.62%  ' }
    @{ Cell = 'E32'; Value = '  -1.92%  ' }
    @{ Cell = 'E33'; Value = '  -1.03%  ' }
    @{ Cell = 'E34'; Value = '  -1.06%  ' }
    @{ Cell = 'E35'; Value = '  +0.11%  ' }
    @{ Cell = 'E36'; Value = '  -3.37%  ' }
    @{ Cell = 'E37'; Value = '  +1.98%  ' }
    @{ Cell = 'E38'; Value = '  -1.13%  ' }
    @{ Cell = 'E39'; Value = '  +0.85%  ' }
    @{ Cell = 'E40'; Value = '  -1.38%  ' }
    @{ Cell = 'E41'; Value = '  -0.82%  ' }
    @{ Cell = 'E42'; Value = '  -0.03%  ' }
    @{ Cell = 'E43'; Value = '  -1.08%  ' }
    @{ Cell = 'E44'; Value = '  +3.29%  ' }
    @{ Cell = 'E45'; Value = '  -1.37%  ' }
    @{ Cell = 'E46'; Value = '  -0.56%  ' }
    @{ Cell = 'E47'; Value = '  -0.09%  ' }
    @{ Cell = 'E48'; Value = '  -1.54%  ' }
    @{ Cell = 'E49'; Value = '  -1.12%  ' }
    @{ Cell = 'E50'; Value = '  -0.16%  ' }
    @{ Cell = 'E51'; Value = '  -0.08%  ' }
)
foreach ($item in $plainCells) {
    $ws.Range($item.Cell).Value = $item.Value
}
